$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.019.63'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '1.919.98'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.90'
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4598'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3826'
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07733'
$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9814'
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.26'
$ws.Range("E11").Value = '  +1.42%  '

$ws.Range("D12").Value = '1.894.13'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.967'
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.689'
$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07005'
$ws.Range("E15").Value = '  -0.79%  '

$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.18'
$ws.Range("E17").Value = '  +1.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009499'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.71'
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").Value = '28.991.42'
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.344'
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.093'
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.51'
$ws.Range("E25").Value = '  +1.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.10'
$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.692'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.94'
$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.855'
$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09335'
$ws.Range("E30").Value = '  +0.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8648'
$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.112'
$ws.Range("E32").Value = '  +0.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.256'
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.037'
$ws.Range("E34").Value = '  -1.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05708'
$ws.Range("E35").Value = '  +1.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02047'
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.522'
$ws.Range("E39").Value = '  +1.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.019'
$ws.Range("E40").Value = '  +12.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5520'
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1751'
$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000002972'
$ws.Range("E43").Value = '  +3.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.386'
$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.217'
$ws.Range("E45").Value = '  +7.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5179'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06902'
$ws.Range("E48").Value = '  +2.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.781'
$ws.Range("E49").Value = '  +0.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.41'

$ws.Range("E51").Value = '  -0.46%  '
